# Watch list test case fixes
# Set the Runmode column (C) to "Y" for all data rows (2-7) on the
# "Test Suite" sheet, and update the selection to reflect the edited range.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

# Update Runmode values to "Y" for rows 2 through 7 (row 3 is already "Y").
$ws.Range("C2:C7").Value = "Y"

# Reflect the edited range as the current selection.
$ws.Activate()
$ws.Range("C2:C7").Select()
